$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$NL = [char]10

# --- Column D width update ---
$ws.Columns.Item(4).ColumnWidth = 66.7109375

# --- Sheet view: zoom / pane / selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true

# --- D column (cross-coverage names) updates ---
$ws.Range("D17").Value = "Cross: FIFO_cov_grp::wr_ack_almostfull_wr_en_cr"
$ws.Range("D18").Value = "Cross: FIFO_cov_grp::wr_ack_empty_wr_en_cr"
$ws.Range("D19").Value = "Cross: FIFO_cov_grp::wr_ack_almostempty_wr_en_cr"
$ws.Range("D20").Value = "Cross: FIFO_cov_grp::wr_ack_full_wr_en_cr" + $NL + "Cross: FIFO_cov_grp::full_almostfull_cr" + $NL + "Cross: FIFO_cov_grp::overflow_cr"
$ws.Range("D22").Value = "Cross: FIFO_cov_grp::almostfull_rd_en_cr"
$ws.Range("D23").Value = "Cross: FIFO_cov_grp::underflow_cr" + $NL + "Cross: FIFO_cov_grp::empty_almostempty_cr"
$ws.Range("D24").Value = "Cross: FIFO_cov_grp::almostempty_rd_en_cr"
$ws.Range("D25").Value = "Cross: FIFO_cov_grp::full_rd_en_cr"

# --- Status column (Pending -> Passed) ---
$ws.Range("I17").Value = "Passed"
$ws.Range("I18").Value = "Passed"
$ws.Range("I19").Value = "Passed"
$ws.Range("I20").Value = "Passed"
$ws.Range("I22").Value = "Passed"
$ws.Range("I23").Value = "Passed"
$ws.Range("I24").Value = "Passed"
$ws.Range("I25").Value = "Passed"

# --- Row heights ---
$ws.Rows.Item(20).RowHeight = 54
$ws.Rows.Item(23).RowHeight = 36
$ws.Rows.Item(26).RowHeight = 72

# --- Restore selection to reported cell ---
$ws.Range("D28").Select()
